# Refresh Leve-profit calculations with the latest Universalis market-board
# price snapshot (currentAveragePrice*, LevePrice*, LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1174.7368
$ws.Range("J17").Value = 1174.7368
$ws.Range("L17").Value = 3524.2104
$ws.Range("N17").Value = -3860.2104

$ws.Range("H43").Value = 2422.3333
$ws.Range("J43").Value = 2393
$ws.Range("L43").Value = 2393
$ws.Range("N43").Value = -2531

$ws.Range("H46").Value = 288428.28
$ws.Range("I46").Value = 334666.34
$ws.Range("J46").Value = 253749.75
$ws.Range("K46").Value = 1003999.02
$ws.Range("L46").Value = 761249.25
$ws.Range("M46").Value = -1003880.02
$ws.Range("N46").Value = -761487.25

$ws.Range("H51").Value = 5000.143
$ws.Range("I51").Value = 5003.5
$ws.Range("J51").Value = 4980
$ws.Range("K51").Value = 5003.5
$ws.Range("L51").Value = 4980
$ws.Range("M51").Value = -4519.5
$ws.Range("N51").Value = -5948

$ws.Range("H60").Value = 288428.28
$ws.Range("I60").Value = 334666.34
$ws.Range("J60").Value = 253749.75
$ws.Range("K60").Value = 1003999.02
$ws.Range("L60").Value = 761249.25
$ws.Range("M60").Value = -1003515.02
$ws.Range("N60").Value = -762217.25

$ws.Range("H98").Value = 53040.316
$ws.Range("I98").Value = 70061.07000000001
$ws.Range("J98").Value = 5382.2
$ws.Range("K98").Value = 70061.07000000001
$ws.Range("L98").Value = 5382.2
$ws.Range("M98").Value = -68563.07000000001
$ws.Range("N98").Value = -8378.200000000001

$ws.Range("H103").Value = 872.5
$ws.Range("I103").Value = 355.58334
$ws.Range("J103").Value = 1492.8
$ws.Range("K103").Value = 1066.75002
$ws.Range("L103").Value = 4478.4
$ws.Range("M103").Value = -480.7500199999999
$ws.Range("N103").Value = -5650.4

$ws.Range("H122").Value = 53040.316
$ws.Range("I122").Value = 70061.07000000001
$ws.Range("J122").Value = 5382.2
$ws.Range("K122").Value = 210183.21
$ws.Range("L122").Value = 16146.6
$ws.Range("M122").Value = -207733.21
$ws.Range("N122").Value = -21046.6

$ws.Range("H141").Value = 5911.9414
$ws.Range("I141").Value = 4121.963
$ws.Range("K141").Value = 12365.889
$ws.Range("M141").Value = -7185.888999999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 43247
$ws.Range("I34").Value = 49994
$ws.Range("J34").Value = 36500
$ws.Range("K34").Value = 49994
$ws.Range("L34").Value = 36500
$ws.Range("M34").Value = -49723
$ws.Range("N34").Value = -37042

$ws.Range("H74").Value = 423295.3
$ws.Range("I74").Value = 682731.6
$ws.Range("J74").Value = 77380.25
$ws.Range("K74").Value = 682731.6
$ws.Range("L74").Value = 77380.25
$ws.Range("M74").Value = -681857.6
$ws.Range("N74").Value = -79128.25

$ws.Range("H77").Value = 423295.3
$ws.Range("I77").Value = 682731.6
$ws.Range("J77").Value = 77380.25
$ws.Range("K77").Value = 3413658
$ws.Range("L77").Value = 386901.25
$ws.Range("M77").Value = -3409290
$ws.Range("N77").Value = -395637.25

$ws.Range("H97").Value = 32258846
$ws.Range("I97").Value = 423.69232
$ws.Range("K97").Value = 423.69232
$ws.Range("M97").Value = 72.30768

$ws.Range("H102").Value = 10755511
$ws.Range("I102").Value = 12348272
$ws.Range("J102").Value = 4377.5
$ws.Range("K102").Value = 12348272
$ws.Range("L102").Value = 4377.5
$ws.Range("M102").Value = -12346650
$ws.Range("N102").Value = -7621.5

$ws.Range("H122").Value = 21436.61
$ws.Range("I122").Value = 1690.2646
$ws.Range("J122").Value = 77384.586
$ws.Range("K122").Value = 5070.793799999999
$ws.Range("L122").Value = 232153.758
$ws.Range("M122").Value = -2620.793799999999
$ws.Range("N122").Value = -237053.758

$ws.Range("H132").Value = 2106.1836
$ws.Range("I132").Value = 2080.1282
$ws.Range("J132").Value = 2207.8
$ws.Range("K132").Value = 6240.3846
$ws.Range("L132").Value = 6623.400000000001
$ws.Range("M132").Value = -3710.3846
$ws.Range("N132").Value = -11683.4


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 319.2
$ws.Range("I5").Value = 275.5
$ws.Range("J5").Value = 494
$ws.Range("K5").Value = 275.5
$ws.Range("L5").Value = 494
$ws.Range("M5").Value = -162.5
$ws.Range("N5").Value = -720

$ws.Range("H7").Value = 912.3333
$ws.Range("I7").Value = 867
$ws.Range("J7").Value = 1003
$ws.Range("K7").Value = 867
$ws.Range("L7").Value = 1003
$ws.Range("M7").Value = -754
$ws.Range("N7").Value = -1229

$ws.Range("H20").Value = 1256.0303
$ws.Range("I20").Value = 1293.875
$ws.Range("J20").Value = 1155.1111
$ws.Range("K20").Value = 1293.875
$ws.Range("L20").Value = 1155.1111
$ws.Range("M20").Value = -1046.875
$ws.Range("N20").Value = -1649.1111

$ws.Range("H80").Value = 951.5
$ws.Range("I80").Value = 2770.8333
$ws.Range("J80").Value = 171.78572
$ws.Range("K80").Value = 2770.8333
$ws.Range("L80").Value = 171.78572
$ws.Range("M80").Value = -1772.8333
$ws.Range("N80").Value = -2167.78572

$ws.Range("H83").Value = 951.5
$ws.Range("I83").Value = 2770.8333
$ws.Range("J83").Value = 171.78572
$ws.Range("K83").Value = 13854.1665
$ws.Range("L83").Value = 858.9286
$ws.Range("M83").Value = -8862.166499999999
$ws.Range("N83").Value = -10842.9286

$ws.Range("H94").Value = 30593.566
$ws.Range("I94").Value = 550.9231
$ws.Range("J94").Value = 225870.75
$ws.Range("K94").Value = 550.9231
$ws.Range("L94").Value = 225870.75
$ws.Range("M94").Value = -99.92309999999998
$ws.Range("N94").Value = -226772.75

$ws.Range("H99").Value = 1956.1538
$ws.Range("I99").Value = 2180.9
$ws.Range("K99").Value = 2180.9
$ws.Range("M99").Value = -682.9000000000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1939.4814
$ws.Range("I31").Value = 1593.4884
$ws.Range("J31").Value = 3292
$ws.Range("K31").Value = 1593.4884
$ws.Range("L31").Value = 3292
$ws.Range("M31").Value = -1298.4884
$ws.Range("N31").Value = -3882

$ws.Range("H34").Value = 1939.4814
$ws.Range("I34").Value = 1593.4884
$ws.Range("J34").Value = 3292
$ws.Range("K34").Value = 1593.4884
$ws.Range("L34").Value = 3292
$ws.Range("M34").Value = -1391.4884
$ws.Range("N34").Value = -3696

$ws.Range("H35").Value = 1634.9667
$ws.Range("I35").Value = 1248.1333
$ws.Range("J35").Value = 2021.8
$ws.Range("K35").Value = 1248.1333
$ws.Range("L35").Value = 2021.8
$ws.Range("M35").Value = -954.1333
$ws.Range("N35").Value = -2609.8

$ws.Range("H69").Value = 22961
$ws.Range("I69").Value = 18454.5
$ws.Range("J69").Value = 50000
$ws.Range("K69").Value = 18454.5
$ws.Range("L69").Value = 50000
$ws.Range("M69").Value = -17705.5
$ws.Range("N69").Value = -51498

$ws.Range("H72").Value = 22961
$ws.Range("I72").Value = 18454.5
$ws.Range("J72").Value = 50000
$ws.Range("K72").Value = 55363.5
$ws.Range("L72").Value = 150000
$ws.Range("M72").Value = -51619.5
$ws.Range("N72").Value = -157488

$ws.Range("H93").Value = 26690.615
$ws.Range("I93").Value = 13840.6
$ws.Range("J93").Value = 34721.875
$ws.Range("K93").Value = 13840.6
$ws.Range("L93").Value = 34721.875
$ws.Range("M93").Value = -11968.6
$ws.Range("N93").Value = -38465.875


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1629
$ws.Range("J18").Value = 3204
$ws.Range("L18").Value = 9612
$ws.Range("N18").Value = -9950

$ws.Range("H50").Value = 979.6667
$ws.Range("I50").Value = 976.2
$ws.Range("K50").Value = 2928.6
$ws.Range("M50").Value = -2447.6

$ws.Range("H53").Value = 979.6667
$ws.Range("I53").Value = 976.2
$ws.Range("K53").Value = 2928.6
$ws.Range("M53").Value = -2447.6

$ws.Range("H132").Value = 2741.8572
$ws.Range("I132").Value = 769.4286
$ws.Range("K132").Value = 6924.8574
$ws.Range("M132").Value = -4394.8574


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 23237
$ws.Range("J35").Value = 39000
$ws.Range("L35").Value = 39000
$ws.Range("N35").Value = -39596

$ws.Range("H132").Value = 2346.6
$ws.Range("I132").Value = 2415.5454
$ws.Range("K132").Value = 7246.6362
$ws.Range("M132").Value = -4716.6362


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4481.3706
$ws.Range("I40").Value = 4535.5
$ws.Range("J40").Value = 4326.7144
$ws.Range("K40").Value = 4535.5
$ws.Range("L40").Value = 4326.7144
$ws.Range("M40").Value = -4399.5
$ws.Range("N40").Value = -4598.7144

$ws.Range("H46").Value = 2928.3015
$ws.Range("J46").Value = 2971.5
$ws.Range("L46").Value = 2971.5
$ws.Range("N46").Value = -3347.5

$ws.Range("H68").Value = 2449.5
$ws.Range("I68").Value = 1600
$ws.Range("K68").Value = 1600
$ws.Range("M68").Value = -851

$ws.Range("H71").Value = 2449.5
$ws.Range("I71").Value = 1600
$ws.Range("K71").Value = 8000
$ws.Range("M71").Value = -4256

$ws.Range("H82").Value = 1178.0834
$ws.Range("I82").Value = 1014.5
$ws.Range("K82").Value = 1014.5
$ws.Range("M82").Value = -653.5

$ws.Range("H85").Value = 1178.0834
$ws.Range("I85").Value = 1014.5
$ws.Range("K85").Value = 1014.5
$ws.Range("M85").Value = 233.5

$ws.Range("H87").Value = 108357.14
$ws.Range("J87").Value = 108357.14
$ws.Range("L87").Value = 108357.14
$ws.Range("N87").Value = -110603.14

$ws.Range("H90").Value = 108357.14
$ws.Range("J90").Value = 108357.14
$ws.Range("L90").Value = 325071.42
$ws.Range("N90").Value = -336303.42

$ws.Range("H122").Value = 3072.228
$ws.Range("I122").Value = 3033.611
$ws.Range("K122").Value = 9100.832999999999
$ws.Range("M122").Value = -6650.832999999999

$ws.Range("H132").Value = 4040.8
$ws.Range("I132").Value = 3007.353
$ws.Range("J132").Value = 6236.875
$ws.Range("K132").Value = 9022.059000000001
$ws.Range("L132").Value = 18710.625
$ws.Range("M132").Value = -6492.059000000001
$ws.Range("N132").Value = -23770.625


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 95000
$ws.Range("J123").Value = 95000
$ws.Range("L123").Value = 95000
$ws.Range("N123").Value = -104800

$ws.Range("H132").Value = 2998.7778
$ws.Range("I132").Value = 3008.111
$ws.Range("J132").Value = 2970.7778
$ws.Range("K132").Value = 9024.332999999999
$ws.Range("L132").Value = 8912.3334
$ws.Range("M132").Value = -6494.332999999999
$ws.Range("N132").Value = -13972.3334

